$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Year headers
$ws.Range("A1").Value = 2017
$ws.Range("A13").Value = 2018

# Product names (column A, rows 2-11) and unit price (column B)
# (written first so the shared-strings table lists products before months)
$products = @("Paper","Pencils","Rulers","Markers","Scissors","Binders","Calculators","Composition Books","Pencil Sharpeners","Erasers")
$prices = @(0.99, 3.99, 0.49, 3.99, 3.29, 2.89, 12.99, 2.99, 0.49, 1.49)
for ($i = 0; $i -lt $products.Length; $i++) {
    $row = $i + 2
    $ws.Range("A" + $row).Value = $products[$i]
    $ws.Range("B" + $row).Value = $prices[$i]
}

# Month headers (row 1, columns C:N)
$months = @("January","February","March","April ","May","June","July","August","September","October","November","December")
$cols = @("C","D","E","F","G","H","I","J","K","L","M","N")
for ($i = 0; $i -lt $months.Length; $i++) {
    $ws.Range($cols[$i] + "1").Value = $months[$i]
}

# Monthly sales figures
$ws.Range("C2").Value = 406
$ws.Range("D2").Value = 307
$ws.Range("E2").Value = 397
$ws.Range("F2").Value = 420
$ws.Range("G2").Value = 356
$ws.Range("H2").Value = 298
$ws.Range("I2").Value = 745
$ws.Range("J2").Value = 912
$ws.Range("K2").Value = 623
$ws.Range("L2").Value = 436
$ws.Range("M2").Value = 300
$ws.Range("N2").Value = 354

$ws.Range("C3").Value = 222
$ws.Range("D3").Value = 235
$ws.Range("E3").Value = 245
$ws.Range("F3").Value = 217
$ws.Range("G3").Value = 197
$ws.Range("H3").Value = 97
$ws.Range("I3").Value = 412
$ws.Range("J3").Value = 467
$ws.Range("K3").Value = 304
$ws.Range("L3").Value = 200
$ws.Range("M3").Value = 135
$ws.Range("N3").Value = 224

$ws.Range("C4").Value = 30
$ws.Range("I4").Value = 45
$ws.Range("J4").Value = 60

# Selection to mimic original author's saved cursor position
$ws.Range("L4").Select()
